# EVENT CUSTOM BUTTONS workbook - "Commands" sheet update
#
# Inserts a new "frequency(<value>[,<sn>])" command row into the PHIDGET
# PWM Command block (enables hi-res internal PID output / PWM frequency
# control), pushing all following rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# Insert a new row at position 64 (just after the PWM Command "out(...)"
# row), shifting rows 64-158 down to 65-159.
$ws.Rows.Item(64).Insert()

# The inserted row inherits formatting from the row below in column A;
# this row has no entry in column A (it's a continuation row), so clear it.
$ws.Range("A64").Clear()

# Fill in the new command documentation cells (columns B and C use the
# same styles as the other rows in this block).
$ws.Range("B64").Value = "frequency(<value>[,<sn>])"
$ws.Range("C64").Value = "PHIDGET PWM Frequency: <value> in Hz"

# Restore the view/selection state recorded for this sheet.
$ws.Activate()
$ws.Range("C63").Select()
$excel.ActiveWindow.ScrollRow = 48
$excel.ActiveWindow.ScrollColumn = 2
